$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header C1 from "audioFalse" to "currentPhase"
$ws.Range("C1").Value = "currentPhase"

# C2 and C3 both now reference a shared "train1P2" string instead of
# their previous distinct audio file paths
$ws.Range("C2").Value = "train1P2"
$ws.Range("C3").Value = "train1P2"
